$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.057.63'
$ws.Range("E2").Value = '''  -5.14%  '
$ws.Range("D3").Value = '''2.224.25'
$ws.Range("E3").Value = '''  -5.94%  '
$ws.Range("E4").Value = '''  +0.11%  '
$ws.Range("D5").Value = '''317.95'
$ws.Range("E5").Value = '''  +2.54%  '
$ws.Range("D6").Value = '''99.07'
$ws.Range("E6").Value = '''  -9.49%  '
$ws.Range("D7").Value = '''0.579'
$ws.Range("E7").Value = '''  -7.84%  '
$ws.Range("E8").Value = '''  +0.12%  '
$ws.Range("D9").Value = '''0.565'
$ws.Range("E9").Value = '''  -8.49%  '
$ws.Range("D10").Value = '''36.54'
$ws.Range("E10").Value = '''  -11.42%  '
$ws.Range("D11").Value = '''54.32'
$ws.Range("E11").Value = '''  -2.02%  '
$ws.Range("E12").Value = '''  -10.23%  '
$ws.Range("E13").Value = '''  -9.18%  '
$ws.Range("E14").Value = '''  -3.96%  '
$ws.Range("D15").Value = '''0.864'
$ws.Range("E15").Value = '''  -12.42%  '
$ws.Range("D16").Value = '''2.565.17'
$ws.Range("E16").Value = '''  -5.91%  '
$ws.Range("D17").Value = '''14.05'
$ws.Range("E17").Value = '''  -8.71%  '
$ws.Range("D18").Value = '''2.217.90'
$ws.Range("E18").Value = '''  -7.24%  '
$ws.Range("D19").Value = '''42.954.39'
$ws.Range("E19").Value = '''  -5.37%  '
$ws.Range("D20").Value = '''14.52'
$ws.Range("E20").Value = '''  +6.41%  '
$ws.Range("D21").Value = '''0.0₃0961'
$ws.Range("E21").Value = '''  -9.77%  '
$ws.Range("D22").Value = '''6.45'
$ws.Range("E22").Value = '''  -11.81%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''65.22'
$ws.Range("E23").Value = '''  -11.12%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '''3.17'
$ws.Range("E24").Value = '''  -8.96%  '
$ws.Range("D25").Value = '''235.95'
$ws.Range("E25").Value = '''  -8.89%  '
$ws.Range("D26").Value = '''2.13'
$ws.Range("E26").Value = '''  -8.22%  '
$ws.Range("E27").Value = '''  +0.38%  '
$ws.Range("D28").Value = '''10.16'
$ws.Range("E28").Value = '''  -9.09%  '
$ws.Range("D29").Value = '''2.19'
$ws.Range("E29").Value = '''  -7.39%  '
$ws.Range("D30").Value = '''6.32'
$ws.Range("E30").Value = '''  -14.03%  '
$ws.Range("D31").Value = '''0.0885'
$ws.Range("E31").Value = '''  -8.77%  '
$ws.Range("D32").Value = '''20.49'
$ws.Range("E32").Value = '''  -8.38%  '
$ws.Range("D33").Value = '''157.62'
$ws.Range("E33").Value = '''  -6.95%  '
$ws.Range("D34").Value = '''33.86'
$ws.Range("E34").Value = '''  -10.99%  '
$ws.Range("E35").Value = '''  -5.56%  '
$ws.Range("E36").Value = '''  +11.82%  '
$ws.Range("D37").Value = '''2.02'
$ws.Range("E37").Value = '''  +15.85%  '
$ws.Range("D38").Value = '''0.122'
$ws.Range("E38").Value = '''  -6.39%  '
$ws.Range("D39").Value = '''4.48'
$ws.Range("E39").Value = '''  -7.47%  '
$ws.Range("E40").Value = '''  -10.62%  '
$ws.Range("D41").Value = '''3.65'
$ws.Range("E41").Value = '''  -7.16%  '
$ws.Range("D42").Value = '''0.0325'
$ws.Range("E42").Value = '''  -9.03%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '''1.854.97'
$ws.Range("E43").Value = '''  +10.74%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '''  +0.02%  '
$ws.Range("D45").Value = '''12.16'
$ws.Range("E45").Value = '''  -6.37%  '
$ws.Range("D46").Value = '''87.83'
$ws.Range("E46").Value = '''  -11.30%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").Value = '''5.48'
$ws.Range("E47").Value = '''  -0.19%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '''0.206'
$ws.Range("E48").Value = '''  -11.83%  '
$ws.Range("D49").Value = '''78.11'
$ws.Range("E49").Value = '''  -6.11%  '
$ws.Range("D50").Value = '''60.50'
$ws.Range("E50").Value = '''  -13.50%  '
$ws.Range("D51").Value = '''8.65'
$ws.Range("E51").Value = '''  -5.51%  '
